$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (players reshuffled/updated per upstream diff)
$data = @(
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Dalton Knecht", "SG", "Los Angeles Lakers"),
    @("Anthony Black", "PG,SG", "Orlando Magic"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Jared McCain", "PG,SG", "Philadelphia 76ers"),
    @("Derrick White", "PG,SG", "Boston Celtics")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
